$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 211.75
$ws.Range("I5").Value = 107.14286
$ws.Range("J5").Value = 358.2
$ws.Range("K5").Value = 107.14286
$ws.Range("L5").Value = 358.2
$ws.Range("M5").Value = 7.857140000000001
$ws.Range("N5").Value = -588.2

$ws.Range("H9").Value = 869.8461
$ws.Range("I9").Value = 254.83333
$ws.Range("J9").Value = 1397
$ws.Range("K9").Value = 254.83333
$ws.Range("L9").Value = 1397
$ws.Range("M9").Value = -85.83332999999999
$ws.Range("N9").Value = -1735

$ws.Range("H70").Value = 2190.3076
$ws.Range("I70").Value = 1698.4
$ws.Range("J70").Value = 2497.75
$ws.Range("K70").Value = 5095.200000000001
$ws.Range("L70").Value = 7493.25
$ws.Range("M70").Value = -4825.200000000001
$ws.Range("N70").Value = -8033.25

$ws.Range("H73").Value = 2190.3076
$ws.Range("I73").Value = 1698.4
$ws.Range("J73").Value = 2497.75
$ws.Range("K73").Value = 5095.200000000001
$ws.Range("L73").Value = 7493.25
$ws.Range("M73").Value = -4159.200000000001
$ws.Range("N73").Value = -9365.25

$ws.Range("H86").Value = 3228.2104
$ws.Range("I86").Value = 1879.9
$ws.Range("J86").Value = 4726.3335
$ws.Range("K86").Value = 1879.9
$ws.Range("L86").Value = 4726.3335
$ws.Range("M86").Value = -756.9000000000001
$ws.Range("N86").Value = -6972.3335

$ws.Range("H89").Value = 3228.2104
$ws.Range("I89").Value = 1879.9
$ws.Range("J89").Value = 4726.3335
$ws.Range("K89").Value = 9399.5
$ws.Range("L89").Value = 23631.6675
$ws.Range("M89").Value = -3783.5
$ws.Range("N89").Value = -34863.6675

$ws.Range("H98").Value = 1057.3334
$ws.Range("I98").Value = 1057.3334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1057.3334
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 440.6666
$ws.Range("N98").ClearContents()

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H122").Value = 1057.3334
$ws.Range("I122").Value = 1057.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3172.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -722.0001999999999
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2276.389
$ws.Range("I132").Value = 1931.6666
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5794.9998
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3264.9998
$ws.Range("N132").Value = -17060

$ws.Range("H135").Value = 1834.5834
$ws.Range("I135").Value = 278.57144
$ws.Range("J135").Value = 4013
$ws.Range("K135").Value = 2507.14296
$ws.Range("L135").Value = 36117
$ws.Range("M135").Value = 27.85703999999987
$ws.Range("N135").Value = -41187

$ws.Range("H138").Value = 3026
$ws.Range("I138").Value = 2561.625
$ws.Range("J138").Value = 3194.8635
$ws.Range("K138").Value = 7684.875
$ws.Range("L138").Value = 9584.5905
$ws.Range("M138").Value = -2544.875
$ws.Range("N138").Value = -19864.5905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2293.5173
$ws.Range("I32").Value = 2095.1482
$ws.Range("J32").Value = 4971.5
$ws.Range("K32").Value = 2095.1482
$ws.Range("L32").Value = 4971.5
$ws.Range("M32").Value = -1808.1482
$ws.Range("N32").Value = -5545.5

$ws.Range("H35").Value = 10572.75
$ws.Range("I35").Value = 2125
$ws.Range("J35").Value = 19020.5
$ws.Range("K35").Value = 2125
$ws.Range("L35").Value = 19020.5
$ws.Range("M35").Value = -1719
$ws.Range("N35").Value = -19832.5

$ws.Range("H97").Value = 3670.3333
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 4005.5
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 4005.5
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -4997.5

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 622.8182
$ws.Range("I20").Value = 423.875
$ws.Range("J20").Value = 1153.3334
$ws.Range("K20").Value = 423.875
$ws.Range("L20").Value = 1153.3334
$ws.Range("M20").Value = -176.875
$ws.Range("N20").Value = -1647.3334

$ws.Range("H86").Value = 2405.2222
$ws.Range("I86").Value = 3019.8
$ws.Range("J86").Value = 1637
$ws.Range("K86").Value = 3019.8
$ws.Range("L86").Value = 1637
$ws.Range("M86").Value = -1896.8
$ws.Range("N86").Value = -3883

$ws.Range("H89").Value = 2405.2222
$ws.Range("I89").Value = 3019.8
$ws.Range("J89").Value = 1637
$ws.Range("K89").Value = 15099
$ws.Range("L89").Value = 8185
$ws.Range("M89").Value = -9483
$ws.Range("N89").Value = -19417

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1638.2
$ws.Range("I58").Value = 1638.2
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1638.2
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1435.2
$ws.Range("N58").ClearContents()

$ws.Range("H86").Value = 12316927
$ws.Range("I86").Value = 14520353
$ws.Range("J86").Value = 1299799
$ws.Range("K86").Value = 14520353
$ws.Range("L86").Value = 1299799
$ws.Range("M86").Value = -14519230
$ws.Range("N86").Value = -1302045

$ws.Range("H89").Value = 12316927
$ws.Range("I89").Value = 14520353
$ws.Range("J89").Value = 1299799
$ws.Range("K89").Value = 72601765
$ws.Range("L89").Value = 6498995
$ws.Range("M89").Value = -72596149
$ws.Range("N89").Value = -6510227

$ws.Range("H99").Value = 3059.7
$ws.Range("I99").Value = 3849.5
$ws.Range("J99").Value = 1875
$ws.Range("K99").Value = 3849.5
$ws.Range("L99").Value = 1875
$ws.Range("M99").Value = -2351.5
$ws.Range("N99").Value = -4871

$ws.Range("H107").Value = 360.42856
$ws.Range("I107").Value = 344.875
$ws.Range("J107").Value = 381.16666
$ws.Range("K107").Value = 344.875
$ws.Range("L107").Value = 381.16666
$ws.Range("M107").Value = 1575.125
$ws.Range("N107").Value = -4221.16666

$ws.Range("H126").Value = 3059.7
$ws.Range("I126").Value = 3849.5
$ws.Range("J126").Value = 1875
$ws.Range("K126").Value = 11548.5
$ws.Range("L126").Value = 5625
$ws.Range("M126").Value = -9078.5
$ws.Range("N126").Value = -10565

$ws.Range("H132").Value = 4371.25
$ws.Range("I132").Value = 4424.2856
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 13272.8568
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -10742.8568
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 1638.2
$ws.Range("I136").Value = 1638.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4914.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2364.6
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1289.5454
$ws.Range("I18").Value = 457
$ws.Range("J18").Value = 1983.3334
$ws.Range("K18").Value = 1371
$ws.Range("L18").Value = 5950.0002
$ws.Range("M18").Value = -1202
$ws.Range("N18").Value = -6288.0002

$ws.Range("H118").Value = 1300
$ws.Range("I118").Value = 1300
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 3900
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -2657

$ws.Range("H132").Value = 2309.077
$ws.Range("I132").Value = 2168.5557
$ws.Range("J132").Value = 2625.25
$ws.Range("K132").Value = 19517.0013
$ws.Range("L132").Value = 23627.25
$ws.Range("M132").Value = -16987.0013
$ws.Range("N132").Value = -28687.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9966.666999999999
$ws.Range("I70").Value = 9966.666999999999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9966.666999999999
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9696.666999999999

$ws.Range("H73").Value = 9966.666999999999
$ws.Range("I73").Value = 9966.666999999999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9966.666999999999
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -9030.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1416.5
$ws.Range("I35").Value = 1416.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1416.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1080.5

$ws.Range("H40").Value = 2949.6667
$ws.Range("I40").Value = 2949.6667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2949.6667
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2813.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1271.8889
$ws.Range("I136").Value = 1271.8889
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3815.6667
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1265.6667
$ws.Range("N136").ClearContents()
